$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 23810702
$ws.Range("J112").Value = 1252.1842
$ws.Range("L112").Value = 3756.5526
$ws.Range("N112").Value = -5972.5526
$ws.Range("H116").Value = 363739.9
$ws.Range("I116").Value = 772085
$ws.Range("J116").Value = 9840.799999999999
$ws.Range("K116").Value = 772085
$ws.Range("L116").Value = 9840.799999999999
$ws.Range("M116").Value = -768643
$ws.Range("N116").Value = -16724.8
$ws.Range("H121").Value = 3000
$ws.Range("J121").Value = 3000
$ws.Range("L121").Value = 9000
$ws.Range("N121").Value = -12494
$ws.Range("H129").Value = 853.51
$ws.Range("I129").Value = 347.75
$ws.Range("J129").Value = 897.48914
$ws.Range("K129").Value = 1043.25
$ws.Range("L129").Value = 2692.46742
$ws.Range("M129").Value = 3956.75
$ws.Range("N129").Value = -12692.46742
$ws.Range("H132").Value = 30610374
$ws.Range("I132").Value = 33337560
$ws.Range("J132").Value = 3338502
$ws.Range("K132").Value = 100012680
$ws.Range("L132").Value = 10015506
$ws.Range("M132").Value = -100010150
$ws.Range("N132").Value = -10020566
$ws.Range("H137").Value = 1192491.1
$ws.Range("I137").Value = 2382526.5
$ws.Range("J137").Value = 2455.75
$ws.Range("K137").Value = 7147579.5
$ws.Range("L137").Value = 7367.25
$ws.Range("M137").Value = -7145029.5
$ws.Range("N137").Value = -12467.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1001.0769
$ws.Range("I2").Value = 820.36365
$ws.Range("K2").Value = 820.36365
$ws.Range("M2").Value = -707.36365
$ws.Range("H32").Value = 3937.6272
$ws.Range("I32").Value = 3865.0176
$ws.Range("J32").Value = 6007
$ws.Range("K32").Value = 3865.0176
$ws.Range("L32").Value = 6007
$ws.Range("M32").Value = -3578.0176
$ws.Range("N32").Value = -6581
$ws.Range("H61").Value = 3405.8372
$ws.Range("I61").Value = 1019.2414
$ws.Range("J61").Value = 8349.5
$ws.Range("K61").Value = 1019.2414
$ws.Range("L61").Value = 8349.5
$ws.Range("M61").Value = -807.2414
$ws.Range("N61").Value = -8773.5
$ws.Range("H74").Value = 2971.5745
$ws.Range("I74").Value = 3202.3547
$ws.Range("J74").Value = 2524.4375
$ws.Range("K74").Value = 3202.3547
$ws.Range("L74").Value = 2524.4375
$ws.Range("M74").Value = -2328.3547
$ws.Range("N74").Value = -4272.4375
$ws.Range("H77").Value = 2971.5745
$ws.Range("I77").Value = 3202.3547
$ws.Range("J77").Value = 2524.4375
$ws.Range("K77").Value = 16011.7735
$ws.Range("L77").Value = 12622.1875
$ws.Range("M77").Value = -11643.7735
$ws.Range("N77").Value = -21358.1875
$ws.Range("H104").Value = 34500
$ws.Range("J104").Value = 34500
$ws.Range("L104").Value = 34500
$ws.Range("N104").Value = -41488
$ws.Range("H110").Value = 1456.5714
$ws.Range("I110").Value = 1333.0588
$ws.Range("K110").Value = 1333.0588
$ws.Range("M110").Value = 711.9412
$ws.Range("H116").Value = 1001.0769
$ws.Range("I116").Value = 820.36365
$ws.Range("K116").Value = 820.36365
$ws.Range("M116").Value = 1473.63635
$ws.Range("H132").Value = 2009.7872
$ws.Range("I132").Value = 1366.6154
$ws.Range("J132").Value = 2806.0952
$ws.Range("K132").Value = 4099.8462
$ws.Range("L132").Value = 8418.285600000001
$ws.Range("M132").Value = -1569.8462
$ws.Range("N132").Value = -13478.2856
$ws.Range("H136").Value = 3405.8372
$ws.Range("I136").Value = 1019.2414
$ws.Range("J136").Value = 8349.5
$ws.Range("K136").Value = 3057.7242
$ws.Range("L136").Value = 25048.5
$ws.Range("M136").Value = -507.7242000000001
$ws.Range("N136").Value = -30148.5
$ws.Range("H137").Value = 45060
$ws.Range("J137").Value = 45060
$ws.Range("L137").Value = 45060
$ws.Range("N137").Value = -55260

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1001.0769
$ws.Range("I3").Value = 820.36365
$ws.Range("K3").Value = 820.36365
$ws.Range("M3").Value = -706.36365
$ws.Range("H80").Value = 174.64516
$ws.Range("I80").Value = 65.90000000000001
$ws.Range("J80").Value = 226.42857
$ws.Range("K80").Value = 65.90000000000001
$ws.Range("L80").Value = 226.42857
$ws.Range("M80").Value = 932.1
$ws.Range("N80").Value = -2222.42857
$ws.Range("H83").Value = 174.64516
$ws.Range("I83").Value = 65.90000000000001
$ws.Range("J83").Value = 226.42857
$ws.Range("K83").Value = 329.5
$ws.Range("L83").Value = 1132.14285
$ws.Range("M83").Value = 4662.5
$ws.Range("N83").Value = -11116.14285
$ws.Range("H105").Value = 1661.5714
$ws.Range("I105").Value = 1664.5333
$ws.Range("J105").Value = 1602.3334
$ws.Range("K105").Value = 1664.5333
$ws.Range("L105").Value = 1602.3334
$ws.Range("M105").Value = 82.46669999999995
$ws.Range("N105").Value = -5096.3334
$ws.Range("H137").Value = 45720
$ws.Range("J137").Value = 45720
$ws.Range("L137").Value = 45720
$ws.Range("N137").Value = -55920

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2781.7144
$ws.Range("I31").Value = 1143.7
$ws.Range("K31").Value = 1143.7
$ws.Range("M31").Value = -848.7
$ws.Range("H34").Value = 2781.7144
$ws.Range("I34").Value = 1143.7
$ws.Range("K34").Value = 1143.7
$ws.Range("M34").Value = -941.7
$ws.Range("H58").Value = 2490.2222
$ws.Range("I58").Value = 1614.5238
$ws.Range("K58").Value = 1614.5238
$ws.Range("M58").Value = -1411.5238
$ws.Range("H81").Value = 26000
$ws.Range("J81").Value = 26000
$ws.Range("L81").Value = 26000
$ws.Range("N81").Value = -27996
$ws.Range("H84").Value = 26000
$ws.Range("J84").Value = 26000
$ws.Range("L84").Value = 78000
$ws.Range("N84").Value = -87984
$ws.Range("H94").Value = 1144.8846
$ws.Range("I94").Value = 975.0833
$ws.Range("J94").Value = 1290.4286
$ws.Range("K94").Value = 975.0833
$ws.Range("L94").Value = 1290.4286
$ws.Range("M94").Value = -524.0833
$ws.Range("N94").Value = -2192.4286
$ws.Range("H105").Value = 1566.9524
$ws.Range("I105").Value = 1329.7858
$ws.Range("J105").Value = 2041.2858
$ws.Range("K105").Value = 1329.7858
$ws.Range("L105").Value = 2041.2858
$ws.Range("M105").Value = 417.2141999999999
$ws.Range("N105").Value = -5535.2858
$ws.Range("H132").Value = 2902.1428
$ws.Range("I132").Value = 2385.4333
$ws.Range("J132").Value = 6002.4
$ws.Range("K132").Value = 7156.2999
$ws.Range("L132").Value = 18007.2
$ws.Range("M132").Value = -4626.2999
$ws.Range("N132").Value = -23067.2
$ws.Range("H136").Value = 2490.2222
$ws.Range("I136").Value = 1614.5238
$ws.Range("K136").Value = 4843.5714
$ws.Range("M136").Value = -2293.5714

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 615.61536
$ws.Range("I113").Value = 602.5
$ws.Range("J113").Value = 626.8570999999999
$ws.Range("K113").Value = 1807.5
$ws.Range("L113").Value = 1880.5713
$ws.Range("M113").Value = 362.5
$ws.Range("N113").Value = -6220.5713
$ws.Range("H131").Value = 784.4
$ws.Range("I131").Value = 355
$ws.Range("J131").Value = 802.2917
$ws.Range("K131").Value = 1065
$ws.Range("L131").Value = 2406.8751
$ws.Range("M131").Value = 3975
$ws.Range("N131").Value = -12486.8751
$ws.Range("H132").Value = 2509.1428
$ws.Range("I132").Value = 519.75
$ws.Range("J132").Value = 3304.9
$ws.Range("K132").Value = 4677.75
$ws.Range("L132").Value = 29744.1
$ws.Range("M132").Value = -2147.75
$ws.Range("N132").Value = -34804.10000000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 10437
$ws.Range("I5").Value = 200
$ws.Range("J5").Value = 12996.25
$ws.Range("K5").Value = 200
$ws.Range("L5").Value = 12996.25
$ws.Range("M5").Value = -88
$ws.Range("N5").Value = -13220.25
$ws.Range("H46").Value = 30679.5
$ws.Range("J46").Value = 30679.5
$ws.Range("L46").Value = 30679.5
$ws.Range("N46").Value = -30991.5
$ws.Range("H126").Value = 3543.33
$ws.Range("I126").Value = 3004.625
$ws.Range("J126").Value = 4928.5713
$ws.Range("K126").Value = 9013.875
$ws.Range("L126").Value = 14785.7139
$ws.Range("M126").Value = -6543.875
$ws.Range("N126").Value = -19725.7139
$ws.Range("H132").Value = 3471.0667
$ws.Range("I132").Value = 3332.5
$ws.Range("J132").Value = 3492.3845
$ws.Range("K132").Value = 9997.5
$ws.Range("L132").Value = 10477.1535
$ws.Range("M132").Value = -7467.5
$ws.Range("N132").Value = -15537.1535
$ws.Range("H137").Value = 40850
$ws.Range("J137").Value = 40850
$ws.Range("L137").Value = 40850
$ws.Range("N137").Value = -51050

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 5849.143
$ws.Range("I122").Value = 2472.5715
$ws.Range("J122").Value = 9225.714
$ws.Range("K122").Value = 7417.7145
$ws.Range("L122").Value = 27677.142
$ws.Range("M122").Value = -4967.7145
$ws.Range("N122").Value = -32577.142

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 432.5
$ws.Range("I107").Value = 400
$ws.Range("J107").Value = 465
$ws.Range("K107").Value = 1200
$ws.Range("L107").Value = 1395
$ws.Range("M107").Value = 720
$ws.Range("N107").Value = -5235
$ws.Range("H132").Value = 2418.84
$ws.Range("I132").Value = 1498.6316
$ws.Range("J132").Value = 5332.8335
$ws.Range("K132").Value = 4495.8948
$ws.Range("L132").Value = 15998.5005
$ws.Range("M132").Value = -1965.8948
$ws.Range("N132").Value = -21058.5005
